$d = $word.ActiveDocument

# The paragraph currently reads "Version 2." split across runs:
#   "Versi" | "on" | " 2" | (bookmark) | "."
# It must become "Version 1." split as:
#   "Version" | " 1." | (bookmark)

# Locate the "Version" text (spans the "Versi" + "on" runs).
$rVersion = $d.Content
$rVersion.Find.Execute("Version", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$versionStart = $rVersion.Start
$versionEnd = $rVersion.End

# Delete the second run ("on") so that rewriting the first run ("Versi")
# below is a genuine text change, which makes Word merge the remaining
# text into a single run instead of leaving the old run split in two.
$rOn = $d.Range($versionStart + 5, $versionEnd)
$rOn.Delete()

# Rewrite "Versi" as the full word "Version".
$rVersi = $d.Range($versionStart, $versionStart + 5)
$rVersi.Text = "Version"

# Change " 2" to " 1." (still one run, now containing the trailing period).
$rNum = $d.Range($versionStart + 7, $versionStart + 9)
$rNum.Text = " 1."

# Remove the old trailing "." run that followed the bookmark.
$rTrailingDot = $d.Range($versionStart + 10, $versionStart + 11)
$rTrailingDot.Delete()
